$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 138
$ws1.Range("F5").Value = 2915
$ws1.Range("F6").Value = 292
$ws1.Range("F7").Value = 396

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 138
$ws4.Range("F5").Value = 2915
$ws4.Range("F6").Value = 292
$ws4.Range("F9").Value = 396
